# Scheduled-runner style update of the Zeromus_Profits leve-crafting data.
# Refreshes currentAveragePrice(NQ/HQ) + derived Leve price/profit columns
# (H:N) across the eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with newly pulled market-board figures. Only value cells are touched;
# a few rows gain/lose their M (LeveProfitNQ) or N (LeveProfitHQ) cell
# entirely depending on whether that recipe now has an NQ/HQ variant.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H104").Value = 3434.3333
$ws.Range("I104").Value = 3434.3333
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 10302.9999
$ws.Range("L104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -8555.999899999999

$ws.Range("H116").Value = 2317768
$ws.Range("I116").Value = 2876439.8
$ws.Range("J116").Value = 3271.2856
$ws.Range("K116").Value = 2876439.8
$ws.Range("L116").Value = 3271.2856
$ws.Range("M116").Value = -2872997.8
$ws.Range("N116").Value = -10155.2856

$ws.Range("H125").Value = 1687.2
$ws.Range("I125").Value = 900
$ws.Range("J125").Value = 1884
$ws.Range("K125").Value = 8100
$ws.Range("L125").Value = 16956
$ws.Range("M125").Value = -5640
$ws.Range("N125").Value = -21876

$ws.Range("H134").Value = 58557
$ws.Range("J134").Value = 58557
$ws.Range("L134").Value = 58557
$ws.Range("N134").Value = -68697

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2224.7446
$ws.Range("I122").Value = 1960.5834
$ws.Range("J122").Value = 3089.2727
$ws.Range("K122").Value = 5881.7502
$ws.Range("L122").Value = 9267.8181
$ws.Range("M122").Value = -3431.7502
$ws.Range("N122").Value = -14167.8181

$ws.Range("H132").Value = 17100.246
$ws.Range("I132").Value = 25163.143
$ws.Range("J132").Value = 2376.6956
$ws.Range("K132").Value = 75489.429
$ws.Range("L132").Value = 7130.0868
$ws.Range("M132").Value = -72959.429
$ws.Range("N132").Value = -12190.0868

$ws.Range("H141").Value = 37637.4
$ws.Range("J141").Value = 37637.4
$ws.Range("L141").Value = 37637.4
$ws.Range("N141").Value = -47997.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 158
$ws.Range("I80").Value = 36
$ws.Range("J80").Value = 166.71428
$ws.Range("K80").Value = 36
$ws.Range("L80").Value = 166.71428
$ws.Range("M80").Value = 962
$ws.Range("N80").Value = -2162.71428

$ws.Range("H83").Value = 158
$ws.Range("I83").Value = 36
$ws.Range("J83").Value = 166.71428
$ws.Range("K83").Value = 180
$ws.Range("L83").Value = 833.5714
$ws.Range("M83").Value = 4812
$ws.Range("N83").Value = -10817.5714

$ws.Range("H86").Value = 4257476.5
$ws.Range("I86").Value = 6251904
$ws.Range("J86").Value = 2698.6
$ws.Range("K86").Value = 6251904
$ws.Range("L86").Value = 2698.6
$ws.Range("M86").Value = -6250781
$ws.Range("N86").Value = -4944.6

$ws.Range("H89").Value = 4257476.5
$ws.Range("I89").Value = 6251904
$ws.Range("J89").Value = 2698.6
$ws.Range("K89").Value = 31259520
$ws.Range("L89").Value = 13493
$ws.Range("M89").Value = -31253904
$ws.Range("N89").Value = -24725

$ws.Range("H94").Value = 12723.962
$ws.Range("I94").Value = 9319.25
$ws.Range("J94").Value = 15642.286
$ws.Range("K94").Value = 9319.25
$ws.Range("L94").Value = 15642.286
$ws.Range("M94").Value = -8868.25
$ws.Range("N94").Value = -16544.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1121.0625
$ws.Range("I16").Value = 1015.44446
$ws.Range("J16").Value = 1256.8572
$ws.Range("K16").Value = 1015.44446
$ws.Range("L16").Value = 1256.8572
$ws.Range("M16").Value = -728.44446
$ws.Range("N16").Value = -1830.8572

$ws.Range("H106").Value = 38000
$ws.Range("J106").Value = 38000
$ws.Range("L106").Value = 38000
$ws.Range("N106").Value = -40524

$ws.Range("H113").Value = 1121.0625
$ws.Range("I113").Value = 1015.44446
$ws.Range("J113").Value = 1256.8572
$ws.Range("K113").Value = 1015.44446
$ws.Range("L113").Value = 1256.8572
$ws.Range("M113").Value = 1154.55554
$ws.Range("N113").Value = -5596.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 814.9167
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 63805510
$ws.Range("I70").Value = 207354770
$ws.Range("J70").Value = 5844.4443
$ws.Range("K70").Value = 207354770
$ws.Range("L70").Value = 5844.4443
$ws.Range("M70").Value = -207354500
$ws.Range("N70").Value = -6384.4443

$ws.Range("H73").Value = 63805510
$ws.Range("I73").Value = 207354770
$ws.Range("J73").Value = 5844.4443
$ws.Range("K73").Value = 207354770
$ws.Range("L73").Value = 5844.4443
$ws.Range("M73").Value = -207353834
$ws.Range("N73").Value = -7716.4443

$ws.Range("H102").Value = 1218.2307
$ws.Range("I102").Value = 1083.7
$ws.Range("J102").Value = 1666.6666
$ws.Range("K102").Value = 1083.7
$ws.Range("L102").Value = 1666.6666
$ws.Range("M102").Value = 538.3
$ws.Range("N102").Value = -4910.6666

$ws.Range("H126").Value = 2309.2163
$ws.Range("I126").Value = 1591.1
$ws.Range("J126").Value = 3154.0588
$ws.Range("K126").Value = 4773.299999999999
$ws.Range("L126").Value = 9462.1764
$ws.Range("M126").Value = -2303.299999999999
$ws.Range("N126").Value = -14402.1764

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1655.091
$ws.Range("I7").Value = 1262.4706
$ws.Range("K7").Value = 1262.4706
$ws.Range("M7").Value = -1150.4706

$ws.Range("H22").Value = 1754920.5
$ws.Range("I22").Value = 3703847
$ws.Range("J22").Value = 886.7
$ws.Range("K22").Value = 3703847
$ws.Range("L22").Value = 886.7
$ws.Range("M22").Value = -3703552
$ws.Range("N22").Value = -1476.7

$ws.Range("H27").Value = 1754920.5
$ws.Range("I27").Value = 3703847
$ws.Range("J27").Value = 886.7
$ws.Range("K27").Value = 3703847
$ws.Range("L27").Value = 886.7
$ws.Range("M27").Value = -3703740
$ws.Range("N27").Value = -1100.7

$ws.Range("H106").Value = 23962
$ws.Range("I106").Value = 4999
$ws.Range("J106").Value = 33443.5
$ws.Range("K106").Value = 4999
$ws.Range("L106").Value = 33443.5
$ws.Range("M106").Value = -3737
$ws.Range("N106").Value = -35967.5

$ws.Range("H126").Value = 1655.091
$ws.Range("I126").Value = 1262.4706
$ws.Range("K126").Value = 3787.4118
$ws.Range("M126").Value = -1317.4118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 16700.666
$ws.Range("J101").Value = 16700.666
$ws.Range("L101").Value = 16700.666
$ws.Range("N101").Value = -23190.666

$ws.Range("H104").Value = 27000
$ws.Range("J104").Value = 27000
$ws.Range("L104").Value = 27000
$ws.Range("N104").Value = -33988

$ws.Range("H105").Value = 49615
$ws.Range("J105").Value = 49615
$ws.Range("L105").Value = 49615
$ws.Range("N105").Value = -56603
